# Bump the published term version 1.0.0 -> 1.1.0 and refresh the
# accompanying "Date" metadata value, as released in term 1.1.0.
#
# Layout of the "Metadata" sheet (column A = Property, column B = Value):
#   Row 3 -> Version
#   Row 8 -> Date

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B3").Value = "1.1.0"
$ws.Range("B8").Value = "2023-07-10T23:08:03+02:00"
